# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeouts) column values for rows 2..39 (column G), replacing
# the previous "Strike#" figures.
$kValues = @(3, 3, 8, 7, 2, 5, 6, 9, 11, 4, 4, 4, 7, 3, 5, 0, 3, 0, 8, 8, 7, 6, 3, 7, 1, 5, 6, 8, 4, 7, 3, 0, 3, 3, 6, 4, 2, 2)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
